$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.167.47"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.085.24"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.83"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.679"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.45"
$ws.Range("E7").Value = "  +18.40%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.393"
$ws.Range("E9").Value = "  +4.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.70"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  +8.14%  "
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.97"
$ws.Range("E13").Value = "  +5.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.390.48"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.827"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.55"
$ws.Range("E16").Value = "  +8.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.086.83"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.152.89"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.72"
$ws.Range("E19").Value = "  +12.96%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +12.12%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.86"
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.50"
$ws.Range("E22").Value = "  +5.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.54"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +16.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.89"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.43"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.49"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").Value = "  +6.80%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.83"
$ws.Range("E32").Value = "  +7.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0642"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +9.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0922"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.119"
$ws.Range("E38").Value = "  +31.36%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.18"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.29"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.52"
$ws.Range("E45").Value = "  +31.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.53"
$ws.Range("E47").Value = "  +14.77%  "
$ws.Range("E48").Value = "  +10.20%  "
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.314.80"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.00"
$ws.Range("E51").Value = "  -0.26%  "
